$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Replicate row 22 formatting (General / date styles) onto the new rows 23-29
$ws.Range("A22:BJ22").Copy()
$ws.Range("A23:BJ29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2) New header cells for columns BK (BRISBANE) and BL (Grau IH)
$ws.Cells.Item(1, 63).NumberFormat = "@"
$ws.Cells.Item(1, 63).Value = 'BRISBANE'
$ws.Cells.Item(1, 63).NumberFormat = "General"
$ws.Cells.Item(1, 64).NumberFormat = "@"
$ws.Cells.Item(1, 64).Value = 'Grau IH'
$ws.Cells.Item(1, 64).NumberFormat = "General"

# --- Row 23 ---
$ws.Cells.Item(23, 1).NumberFormat = "@"
$ws.Cells.Item(23, 1).Value = 'Segmentectomia o Bisegmentectomia'
$ws.Cells.Item(23, 1).NumberFormat = "General"
$ws.Cells.Item(23, 2).Value = 1737
$ws.Cells.Item(23, 3).NumberFormat = "@"
$ws.Cells.Item(23, 3).Value = '05/04/2018'
$ws.Cells.Item(23, 3).NumberFormat = "General"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '09/06/2021'
$ws.Cells.Item(23, 4).NumberFormat = "General"
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = '01/01/2018'
$ws.Cells.Item(23, 5).NumberFormat = "General"
$ws.Cells.Item(23, 6).NumberFormat = "@"
$ws.Cells.Item(23, 6).Value = 'IV,III'
$ws.Cells.Item(23, 6).NumberFormat = "General"
$ws.Cells.Item(23, 7).Value = 1497
$ws.Cells.Item(23, 8).NumberFormat = "@"
$ws.Cells.Item(23, 8).Value = 'Jaume'
$ws.Cells.Item(23, 8).NumberFormat = "General"
$ws.Cells.Item(23, 9).NumberFormat = "@"
$ws.Cells.Item(23, 9).Value = 'Gual2'
$ws.Cells.Item(23, 9).NumberFormat = "General"
$ws.Cells.Item(23, 10).NumberFormat = "@"
$ws.Cells.Item(23, 10).Value = 'Bosch2'
$ws.Cells.Item(23, 10).NumberFormat = "General"
$ws.Cells.Item(23, 11).Value = 13297134
$ws.Cells.Item(23, 12).NumberFormat = "@"
$ws.Cells.Item(23, 12).Value = 'Si'
$ws.Cells.Item(23, 12).NumberFormat = "General"
$ws.Cells.Item(23, 13).NumberFormat = "@"
$ws.Cells.Item(23, 13).Value = 'Home'
$ws.Cells.Item(23, 13).NumberFormat = "General"
$ws.Cells.Item(23, 14).NumberFormat = "@"
$ws.Cells.Item(23, 14).Value = '49'
$ws.Cells.Item(23, 14).NumberFormat = "General"
$ws.Cells.Item(23, 20).Value = 44119
$ws.Cells.Item(23, 21).NumberFormat = "@"
$ws.Cells.Item(23, 21).Value = 'Resecció Menor (<3 segm)'
$ws.Cells.Item(23, 21).NumberFormat = "General"
$ws.Cells.Item(23, 22).NumberFormat = "@"
$ws.Cells.Item(23, 22).Value = 'segmentectomia 4a'
$ws.Cells.Item(23, 22).NumberFormat = "General"
$ws.Cells.Item(23, 23).NumberFormat = "@"
$ws.Cells.Item(23, 23).Value = 'Oberta'
$ws.Cells.Item(23, 23).NumberFormat = "General"
$ws.Cells.Item(23, 24).NumberFormat = "@"
$ws.Cells.Item(23, 24).Value = 'No'
$ws.Cells.Item(23, 24).NumberFormat = "General"
$ws.Cells.Item(23, 25).NumberFormat = "@"
$ws.Cells.Item(23, 25).Value = 'No'
$ws.Cells.Item(23, 25).NumberFormat = "General"
$ws.Cells.Item(23, 26).NumberFormat = "@"
$ws.Cells.Item(23, 26).Value = 'Impressió R0'
$ws.Cells.Item(23, 26).NumberFormat = "General"
$ws.Cells.Item(23, 29).NumberFormat = "@"
$ws.Cells.Item(23, 29).Value = 'No'
$ws.Cells.Item(23, 29).NumberFormat = "General"
$ws.Cells.Item(23, 30).NumberFormat = "@"
$ws.Cells.Item(23, 30).Value = 'Si'
$ws.Cells.Item(23, 30).NumberFormat = "General"
$ws.Cells.Item(23, 31).NumberFormat = "@"
$ws.Cells.Item(23, 31).Value = 'No'
$ws.Cells.Item(23, 31).NumberFormat = "General"
$ws.Cells.Item(23, 32).NumberFormat = "@"
$ws.Cells.Item(23, 32).Value = 'Si'
$ws.Cells.Item(23, 32).NumberFormat = "General"
$ws.Cells.Item(23, 33).NumberFormat = "@"
$ws.Cells.Item(23, 33).Value = 'IIIa'
$ws.Cells.Item(23, 33).NumberFormat = "General"
$ws.Cells.Item(23, 34).Value = 27.6
$ws.Cells.Item(23, 35).Value = 2
$ws.Cells.Item(23, 36).Value = 3
$ws.Cells.Item(23, 37).Value = 0
$ws.Cells.Item(23, 38).NumberFormat = "@"
$ws.Cells.Item(23, 38).Value = 'Si'
$ws.Cells.Item(23, 38).NumberFormat = "General"
$ws.Cells.Item(23, 39).NumberFormat = "@"
$ws.Cells.Item(23, 39).Value = 'es tracta del marhe de la linea de transecció previa'
$ws.Cells.Item(23, 39).NumberFormat = "General"
$ws.Cells.Item(23, 40).Value = 44522
$ws.Cells.Item(23, 41).NumberFormat = "@"
$ws.Cells.Item(23, 41).Value = 'No'
$ws.Cells.Item(23, 41).NumberFormat = "General"
$ws.Cells.Item(23, 42).NumberFormat = "@"
$ws.Cells.Item(23, 42).Value = 'No'
$ws.Cells.Item(23, 42).NumberFormat = "General"
$ws.Cells.Item(23, 43).NumberFormat = "@"
$ws.Cells.Item(23, 43).Value = 'Viu'
$ws.Cells.Item(23, 43).NumberFormat = "General"
$ws.Cells.Item(23, 44).NumberFormat = "@"
$ws.Cells.Item(23, 44).Value = 'Si'
$ws.Cells.Item(23, 44).NumberFormat = "General"
$ws.Cells.Item(23, 45).NumberFormat = "@"
$ws.Cells.Item(23, 45).Value = 'Si'
$ws.Cells.Item(23, 45).NumberFormat = "General"
$ws.Cells.Item(23, 46).NumberFormat = "@"
$ws.Cells.Item(23, 46).Value = 'No'
$ws.Cells.Item(23, 46).NumberFormat = "General"
$ws.Cells.Item(23, 47).NumberFormat = "@"
$ws.Cells.Item(23, 47).Value = 'No'
$ws.Cells.Item(23, 47).NumberFormat = "General"
$ws.Cells.Item(23, 48).NumberFormat = "@"
$ws.Cells.Item(23, 48).Value = 'No'
$ws.Cells.Item(23, 48).NumberFormat = "General"
$ws.Cells.Item(23, 49).NumberFormat = "@"
$ws.Cells.Item(23, 49).Value = 'No'
$ws.Cells.Item(23, 49).NumberFormat = "General"
$ws.Cells.Item(23, 50).NumberFormat = "@"
$ws.Cells.Item(23, 50).Value = 'Si'
$ws.Cells.Item(23, 50).NumberFormat = "General"
$ws.Cells.Item(23, 51).NumberFormat = "@"
$ws.Cells.Item(23, 51).Value = 'Si'
$ws.Cells.Item(23, 51).NumberFormat = "General"
$ws.Cells.Item(23, 52).NumberFormat = "@"
$ws.Cells.Item(23, 52).Value = 'Si'
$ws.Cells.Item(23, 52).NumberFormat = "General"
$ws.Cells.Item(23, 53).Value = 6
$ws.Cells.Item(23, 54).Value = 44610.47931371528
$ws.Cells.Item(23, 55).NumberFormat = "@"
$ws.Cells.Item(23, 55).Value = 'Falta alguna variable, revisar'
$ws.Cells.Item(23, 55).NumberFormat = "General"
$ws.Cells.Item(23, 56).Value = 43342
$ws.Cells.Item(23, 57).NumberFormat = "@"
$ws.Cells.Item(23, 57).Value = 'No'
$ws.Cells.Item(23, 57).NumberFormat = "General"
$ws.Cells.Item(23, 59).Value = 1
$ws.Cells.Item(23, 63).NumberFormat = "@"
$ws.Cells.Item(23, 63).Value = 'Segmentectomia1a8'
$ws.Cells.Item(23, 63).NumberFormat = "General"

# --- Row 24 ---
$ws.Cells.Item(24, 1).NumberFormat = "@"
$ws.Cells.Item(24, 1).Value = 'Hepatectomia major + resecció contralateral'
$ws.Cells.Item(24, 1).NumberFormat = "General"
$ws.Cells.Item(24, 2).Value = 1632
$ws.Cells.Item(24, 3).NumberFormat = "@"
$ws.Cells.Item(24, 3).Value = '05/04/2018'
$ws.Cells.Item(24, 3).NumberFormat = "General"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '09/06/2021'
$ws.Cells.Item(24, 4).NumberFormat = "General"
$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = '01/01/2018'
$ws.Cells.Item(24, 5).NumberFormat = "General"
$ws.Cells.Item(24, 6).NumberFormat = "@"
$ws.Cells.Item(24, 6).Value = 'IV'
$ws.Cells.Item(24, 6).NumberFormat = "General"
$ws.Cells.Item(24, 7).Value = 1497
$ws.Cells.Item(24, 8).NumberFormat = "@"
$ws.Cells.Item(24, 8).Value = 'JAUME'
$ws.Cells.Item(24, 8).NumberFormat = "General"
$ws.Cells.Item(24, 9).NumberFormat = "@"
$ws.Cells.Item(24, 9).Value = 'GUAL'
$ws.Cells.Item(24, 9).NumberFormat = "General"
$ws.Cells.Item(24, 10).NumberFormat = "@"
$ws.Cells.Item(24, 10).Value = 'BOSCH'
$ws.Cells.Item(24, 10).NumberFormat = "General"
$ws.Cells.Item(24, 11).Value = 13297134
$ws.Cells.Item(24, 12).NumberFormat = "@"
$ws.Cells.Item(24, 12).Value = 'Si'
$ws.Cells.Item(24, 12).NumberFormat = "General"
$ws.Cells.Item(24, 13).NumberFormat = "@"
$ws.Cells.Item(24, 13).Value = 'Home'
$ws.Cells.Item(24, 13).NumberFormat = "General"
$ws.Cells.Item(24, 14).NumberFormat = "@"
$ws.Cells.Item(24, 14).Value = '47'
$ws.Cells.Item(24, 14).NumberFormat = "General"
$ws.Cells.Item(24, 15).NumberFormat = "@"
$ws.Cells.Item(24, 15).Value = '87'
$ws.Cells.Item(24, 15).NumberFormat = "General"
$ws.Cells.Item(24, 16).Value = 170
$ws.Cells.Item(24, 17).Value = 30
$ws.Cells.Item(24, 18).Value = 2
$ws.Cells.Item(24, 19).NumberFormat = "@"
$ws.Cells.Item(24, 19).Value = 'Si'
$ws.Cells.Item(24, 19).NumberFormat = "General"
$ws.Cells.Item(24, 20).NumberFormat = "@"
$ws.Cells.Item(24, 20).Value = '02/15/2019'
$ws.Cells.Item(24, 20).NumberFormat = "General"
$ws.Cells.Item(24, 21).NumberFormat = "@"
$ws.Cells.Item(24, 21).Value = 'Resecció Major (>= 3 segm)'
$ws.Cells.Item(24, 21).NumberFormat = "General"
$ws.Cells.Item(24, 22).NumberFormat = "@"
$ws.Cells.Item(24, 22).Value = 'Hepatectomia derecha + reseccion limitada segmento IV'
$ws.Cells.Item(24, 22).NumberFormat = "General"
$ws.Cells.Item(24, 23).NumberFormat = "@"
$ws.Cells.Item(24, 23).Value = 'Oberta'
$ws.Cells.Item(24, 23).NumberFormat = "General"
$ws.Cells.Item(24, 24).NumberFormat = "@"
$ws.Cells.Item(24, 24).Value = 'No'
$ws.Cells.Item(24, 24).NumberFormat = "General"
$ws.Cells.Item(24, 25).NumberFormat = "@"
$ws.Cells.Item(24, 25).Value = 'No'
$ws.Cells.Item(24, 25).NumberFormat = "General"
$ws.Cells.Item(24, 26).NumberFormat = "@"
$ws.Cells.Item(24, 26).Value = 'Impressió R0'
$ws.Cells.Item(24, 26).NumberFormat = "General"
$ws.Cells.Item(24, 27).Value = 10
$ws.Cells.Item(24, 28).Value = 2
$ws.Cells.Item(24, 29).NumberFormat = "@"
$ws.Cells.Item(24, 29).Value = 'Si'
$ws.Cells.Item(24, 29).NumberFormat = "General"
$ws.Cells.Item(24, 30).NumberFormat = "@"
$ws.Cells.Item(24, 30).Value = 'No'
$ws.Cells.Item(24, 30).NumberFormat = "General"
$ws.Cells.Item(24, 31).NumberFormat = "@"
$ws.Cells.Item(24, 31).Value = 'No'
$ws.Cells.Item(24, 31).NumberFormat = "General"
$ws.Cells.Item(24, 32).NumberFormat = "@"
$ws.Cells.Item(24, 32).Value = 'No'
$ws.Cells.Item(24, 32).NumberFormat = "General"
$ws.Cells.Item(24, 33).NumberFormat = "@"
$ws.Cells.Item(24, 33).Value = '0'
$ws.Cells.Item(24, 33).NumberFormat = "General"
$ws.Cells.Item(24, 34).Value = 0
$ws.Cells.Item(24, 35).Value = 2
$ws.Cells.Item(24, 36).Value = 1.3
$ws.Cells.Item(24, 37).Value = 0.2
$ws.Cells.Item(24, 38).NumberFormat = "@"
$ws.Cells.Item(24, 38).Value = 'No'
$ws.Cells.Item(24, 38).NumberFormat = "General"
# AM24: explicit empty string in source (Excel normalizes "" to blank) - skipped
$ws.Cells.Item(24, 40).Value = 44522
$ws.Cells.Item(24, 41).NumberFormat = "@"
$ws.Cells.Item(24, 41).Value = 'No'
$ws.Cells.Item(24, 41).NumberFormat = "General"
$ws.Cells.Item(24, 42).NumberFormat = "@"
$ws.Cells.Item(24, 42).Value = 'No'
$ws.Cells.Item(24, 42).NumberFormat = "General"
$ws.Cells.Item(24, 43).NumberFormat = "@"
$ws.Cells.Item(24, 43).Value = 'Viu'
$ws.Cells.Item(24, 43).NumberFormat = "General"
$ws.Cells.Item(24, 44).NumberFormat = "@"
$ws.Cells.Item(24, 44).Value = 'No'
$ws.Cells.Item(24, 44).NumberFormat = "General"
$ws.Cells.Item(24, 45).NumberFormat = "@"
$ws.Cells.Item(24, 45).Value = 'No'
$ws.Cells.Item(24, 45).NumberFormat = "General"
$ws.Cells.Item(24, 46).NumberFormat = "@"
$ws.Cells.Item(24, 46).Value = 'No'
$ws.Cells.Item(24, 46).NumberFormat = "General"
$ws.Cells.Item(24, 47).NumberFormat = "@"
$ws.Cells.Item(24, 47).Value = 'No'
$ws.Cells.Item(24, 47).NumberFormat = "General"
$ws.Cells.Item(24, 48).NumberFormat = "@"
$ws.Cells.Item(24, 48).Value = 'No'
$ws.Cells.Item(24, 48).NumberFormat = "General"
# AX24: explicit empty string in source (Excel normalizes "" to blank) - skipped
$ws.Cells.Item(24, 51).NumberFormat = "@"
$ws.Cells.Item(24, 51).Value = 'Si'
$ws.Cells.Item(24, 51).NumberFormat = "General"
$ws.Cells.Item(24, 52).NumberFormat = "@"
$ws.Cells.Item(24, 52).Value = 'Si'
$ws.Cells.Item(24, 52).NumberFormat = "General"
$ws.Cells.Item(24, 53).Value = 10
$ws.Cells.Item(24, 54).Value = 44610.4793574537
$ws.Cells.Item(24, 55).NumberFormat = "@"
$ws.Cells.Item(24, 55).Value = 'unable to complete promise all for CMD data after condition CMD = true, error message: '
$ws.Cells.Item(24, 55).NumberFormat = "General"
$ws.Cells.Item(24, 56).Value = 43344
$ws.Cells.Item(24, 57).NumberFormat = "@"
$ws.Cells.Item(24, 57).Value = 'No'
$ws.Cells.Item(24, 57).NumberFormat = "General"
$ws.Cells.Item(24, 58).NumberFormat = "@"
$ws.Cells.Item(24, 58).Value = 'No'
$ws.Cells.Item(24, 58).NumberFormat = "General"
# BH24: explicit empty string in source (Excel normalizes "" to blank) - skipped
# BI24: explicit empty string in source (Excel normalizes "" to blank) - skipped
$ws.Cells.Item(24, 63).NumberFormat = "@"
$ws.Cells.Item(24, 63).Value = 'hepatectomiaDreta,Segmentectomia1a8'
$ws.Cells.Item(24, 63).NumberFormat = "General"
# BL24: explicit empty string in source (Excel normalizes "" to blank) - skipped

# --- Row 25 ---
$ws.Cells.Item(25, 1).NumberFormat = "@"
$ws.Cells.Item(25, 1).Value = 'Hepatectomia dreta'
$ws.Cells.Item(25, 1).NumberFormat = "General"
$ws.Cells.Item(25, 2).Value = 1695
$ws.Cells.Item(25, 3).NumberFormat = "@"
$ws.Cells.Item(25, 3).Value = '09/07/2018'
$ws.Cells.Item(25, 3).NumberFormat = "General"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '11/02/2019'
$ws.Cells.Item(25, 4).NumberFormat = "General"
$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = '25/06/2018'
$ws.Cells.Item(25, 5).NumberFormat = "General"
$ws.Cells.Item(25, 6).NumberFormat = "@"
$ws.Cells.Item(25, 6).Value = 'VII'
$ws.Cells.Item(25, 6).NumberFormat = "General"
$ws.Cells.Item(25, 8).NumberFormat = "@"
$ws.Cells.Item(25, 8).Value = 'Jordi'
$ws.Cells.Item(25, 8).NumberFormat = "General"
$ws.Cells.Item(25, 9).NumberFormat = "@"
$ws.Cells.Item(25, 9).Value = 'Morillas2'
$ws.Cells.Item(25, 9).NumberFormat = "General"
$ws.Cells.Item(25, 10).NumberFormat = "@"
$ws.Cells.Item(25, 10).Value = 'Esteban2'
$ws.Cells.Item(25, 10).NumberFormat = "General"
$ws.Cells.Item(25, 11).Value = 13296015
$ws.Cells.Item(25, 12).NumberFormat = "@"
$ws.Cells.Item(25, 12).Value = 'No'
$ws.Cells.Item(25, 12).NumberFormat = "General"
$ws.Cells.Item(25, 15).NumberFormat = "@"
$ws.Cells.Item(25, 15).Value = '79'
$ws.Cells.Item(25, 15).NumberFormat = "General"
$ws.Cells.Item(25, 16).Value = 178
$ws.Cells.Item(25, 17).Value = 25
$ws.Cells.Item(25, 18).Value = 3
$ws.Cells.Item(25, 19).NumberFormat = "@"
$ws.Cells.Item(25, 19).Value = 'No'
$ws.Cells.Item(25, 19).NumberFormat = "General"
$ws.Cells.Item(25, 20).Value = 43858
$ws.Cells.Item(25, 21).NumberFormat = "@"
$ws.Cells.Item(25, 21).Value = 'Resecció Major (>= 3 segm)'
$ws.Cells.Item(25, 21).NumberFormat = "General"
$ws.Cells.Item(25, 22).NumberFormat = "@"
$ws.Cells.Item(25, 22).Value = 'hepatectomia dreta'
$ws.Cells.Item(25, 22).NumberFormat = "General"
$ws.Cells.Item(25, 23).NumberFormat = "@"
$ws.Cells.Item(25, 23).Value = 'Oberta'
$ws.Cells.Item(25, 23).NumberFormat = "General"
$ws.Cells.Item(25, 24).NumberFormat = "@"
$ws.Cells.Item(25, 24).Value = 'Si, com a primer temps quirúrgic'
$ws.Cells.Item(25, 24).NumberFormat = "General"
$ws.Cells.Item(25, 25).NumberFormat = "@"
$ws.Cells.Item(25, 25).Value = 'No'
$ws.Cells.Item(25, 25).NumberFormat = "General"
$ws.Cells.Item(25, 26).NumberFormat = "@"
$ws.Cells.Item(25, 26).Value = 'Impressió R1'
$ws.Cells.Item(25, 26).NumberFormat = "General"
$ws.Cells.Item(25, 27).Value = 1
$ws.Cells.Item(25, 28).Value = 3
$ws.Cells.Item(25, 29).NumberFormat = "@"
$ws.Cells.Item(25, 29).Value = 'No'
$ws.Cells.Item(25, 29).NumberFormat = "General"
$ws.Cells.Item(25, 30).NumberFormat = "@"
$ws.Cells.Item(25, 30).Value = 'Si'
$ws.Cells.Item(25, 30).NumberFormat = "General"
$ws.Cells.Item(25, 31).NumberFormat = "@"
$ws.Cells.Item(25, 31).Value = 'No'
$ws.Cells.Item(25, 31).NumberFormat = "General"
$ws.Cells.Item(25, 32).NumberFormat = "@"
$ws.Cells.Item(25, 32).Value = 'Si'
$ws.Cells.Item(25, 32).NumberFormat = "General"
$ws.Cells.Item(25, 33).NumberFormat = "@"
$ws.Cells.Item(25, 33).Value = 'IIIb'
$ws.Cells.Item(25, 33).NumberFormat = "General"
$ws.Cells.Item(25, 34).Value = 61
$ws.Cells.Item(25, 35).Value = 1
$ws.Cells.Item(25, 36).Value = 3
$ws.Cells.Item(25, 37).Value = 0
$ws.Cells.Item(25, 38).NumberFormat = "@"
$ws.Cells.Item(25, 38).Value = 'Si'
$ws.Cells.Item(25, 38).NumberFormat = "General"
$ws.Cells.Item(25, 39).NumberFormat = "@"
$ws.Cells.Item(25, 39).Value = 'ampliacio quirurgica'
$ws.Cells.Item(25, 39).NumberFormat = "General"
$ws.Cells.Item(25, 40).Value = 43983
$ws.Cells.Item(25, 41).NumberFormat = "@"
$ws.Cells.Item(25, 41).Value = 'No'
$ws.Cells.Item(25, 41).NumberFormat = "General"
$ws.Cells.Item(25, 42).NumberFormat = "@"
$ws.Cells.Item(25, 42).Value = 'No'
$ws.Cells.Item(25, 42).NumberFormat = "General"
$ws.Cells.Item(25, 43).NumberFormat = "@"
$ws.Cells.Item(25, 43).Value = 'Viu'
$ws.Cells.Item(25, 43).NumberFormat = "General"
$ws.Cells.Item(25, 44).NumberFormat = "@"
$ws.Cells.Item(25, 44).Value = 'Si'
$ws.Cells.Item(25, 44).NumberFormat = "General"
$ws.Cells.Item(25, 45).NumberFormat = "@"
$ws.Cells.Item(25, 45).Value = 'Si'
$ws.Cells.Item(25, 45).NumberFormat = "General"
$ws.Cells.Item(25, 46).NumberFormat = "@"
$ws.Cells.Item(25, 46).Value = 'No'
$ws.Cells.Item(25, 46).NumberFormat = "General"
$ws.Cells.Item(25, 47).NumberFormat = "@"
$ws.Cells.Item(25, 47).Value = 'No'
$ws.Cells.Item(25, 47).NumberFormat = "General"
$ws.Cells.Item(25, 48).NumberFormat = "@"
$ws.Cells.Item(25, 48).Value = 'Si'
$ws.Cells.Item(25, 48).NumberFormat = "General"
$ws.Cells.Item(25, 49).NumberFormat = "@"
$ws.Cells.Item(25, 49).Value = 'No'
$ws.Cells.Item(25, 49).NumberFormat = "General"
$ws.Cells.Item(25, 50).NumberFormat = "@"
$ws.Cells.Item(25, 50).Value = 'No'
$ws.Cells.Item(25, 50).NumberFormat = "General"
$ws.Cells.Item(25, 51).NumberFormat = "@"
$ws.Cells.Item(25, 51).Value = 'No'
$ws.Cells.Item(25, 51).NumberFormat = "General"
$ws.Cells.Item(25, 52).NumberFormat = "@"
$ws.Cells.Item(25, 52).Value = 'No'
$ws.Cells.Item(25, 52).NumberFormat = "General"
$ws.Cells.Item(25, 53).Value = 55
$ws.Cells.Item(25, 54).Value = 44610.47960366898
$ws.Cells.Item(25, 55).NumberFormat = "@"
$ws.Cells.Item(25, 55).Value = 'Error with MidaAp or NumAp'
$ws.Cells.Item(25, 55).NumberFormat = "General"
$ws.Cells.Item(25, 56).Value = 43446
$ws.Cells.Item(25, 57).NumberFormat = "@"
$ws.Cells.Item(25, 57).Value = 'No'
$ws.Cells.Item(25, 57).NumberFormat = "General"
$ws.Cells.Item(25, 58).NumberFormat = "@"
$ws.Cells.Item(25, 58).Value = 'oclusió, peritonitis fecaoidea'
$ws.Cells.Item(25, 58).NumberFormat = "General"
$ws.Cells.Item(25, 59).Value = 2
$ws.Cells.Item(25, 63).NumberFormat = "@"
$ws.Cells.Item(25, 63).Value = 'hepatectomiaDreta'
$ws.Cells.Item(25, 63).NumberFormat = "General"

# --- Row 26 ---
$ws.Cells.Item(26, 1).NumberFormat = "@"
$ws.Cells.Item(26, 1).Value = 'Segmentectomia o Bisegmentectomia'
$ws.Cells.Item(26, 1).NumberFormat = "General"
$ws.Cells.Item(26, 2).Value = 1737
$ws.Cells.Item(26, 3).NumberFormat = "@"
$ws.Cells.Item(26, 3).Value = '05/04/2018'
$ws.Cells.Item(26, 3).NumberFormat = "General"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '09/06/2021'
$ws.Cells.Item(26, 4).NumberFormat = "General"
$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = '01/01/2018'
$ws.Cells.Item(26, 5).NumberFormat = "General"
$ws.Cells.Item(26, 6).NumberFormat = "@"
$ws.Cells.Item(26, 6).Value = 'IV,III'
$ws.Cells.Item(26, 6).NumberFormat = "General"
$ws.Cells.Item(26, 7).Value = 1497
$ws.Cells.Item(26, 8).NumberFormat = "@"
$ws.Cells.Item(26, 8).Value = 'Jaume'
$ws.Cells.Item(26, 8).NumberFormat = "General"
$ws.Cells.Item(26, 9).NumberFormat = "@"
$ws.Cells.Item(26, 9).Value = 'Gual2'
$ws.Cells.Item(26, 9).NumberFormat = "General"
$ws.Cells.Item(26, 10).NumberFormat = "@"
$ws.Cells.Item(26, 10).Value = 'Bosch2'
$ws.Cells.Item(26, 10).NumberFormat = "General"
$ws.Cells.Item(26, 11).Value = 13297134
$ws.Cells.Item(26, 12).NumberFormat = "@"
$ws.Cells.Item(26, 12).Value = 'Si'
$ws.Cells.Item(26, 12).NumberFormat = "General"
$ws.Cells.Item(26, 13).NumberFormat = "@"
$ws.Cells.Item(26, 13).Value = 'Home'
$ws.Cells.Item(26, 13).NumberFormat = "General"
$ws.Cells.Item(26, 14).NumberFormat = "@"
$ws.Cells.Item(26, 14).Value = '49'
$ws.Cells.Item(26, 14).NumberFormat = "General"
$ws.Cells.Item(26, 20).Value = 44119
$ws.Cells.Item(26, 21).NumberFormat = "@"
$ws.Cells.Item(26, 21).Value = 'Resecció Menor (<3 segm)'
$ws.Cells.Item(26, 21).NumberFormat = "General"
$ws.Cells.Item(26, 22).NumberFormat = "@"
$ws.Cells.Item(26, 22).Value = 'segmentectomia 4a'
$ws.Cells.Item(26, 22).NumberFormat = "General"
$ws.Cells.Item(26, 23).NumberFormat = "@"
$ws.Cells.Item(26, 23).Value = 'Oberta'
$ws.Cells.Item(26, 23).NumberFormat = "General"
$ws.Cells.Item(26, 24).NumberFormat = "@"
$ws.Cells.Item(26, 24).Value = 'No'
$ws.Cells.Item(26, 24).NumberFormat = "General"
$ws.Cells.Item(26, 25).NumberFormat = "@"
$ws.Cells.Item(26, 25).Value = 'No'
$ws.Cells.Item(26, 25).NumberFormat = "General"
$ws.Cells.Item(26, 26).NumberFormat = "@"
$ws.Cells.Item(26, 26).Value = 'Impressió R0'
$ws.Cells.Item(26, 26).NumberFormat = "General"
$ws.Cells.Item(26, 29).NumberFormat = "@"
$ws.Cells.Item(26, 29).Value = 'No'
$ws.Cells.Item(26, 29).NumberFormat = "General"
$ws.Cells.Item(26, 30).NumberFormat = "@"
$ws.Cells.Item(26, 30).Value = 'Si'
$ws.Cells.Item(26, 30).NumberFormat = "General"
$ws.Cells.Item(26, 31).NumberFormat = "@"
$ws.Cells.Item(26, 31).Value = 'No'
$ws.Cells.Item(26, 31).NumberFormat = "General"
$ws.Cells.Item(26, 32).NumberFormat = "@"
$ws.Cells.Item(26, 32).Value = 'Si'
$ws.Cells.Item(26, 32).NumberFormat = "General"
$ws.Cells.Item(26, 33).NumberFormat = "@"
$ws.Cells.Item(26, 33).Value = 'IIIa'
$ws.Cells.Item(26, 33).NumberFormat = "General"
$ws.Cells.Item(26, 34).Value = 27.6
$ws.Cells.Item(26, 35).Value = 2
$ws.Cells.Item(26, 36).Value = 3
$ws.Cells.Item(26, 37).Value = 0
$ws.Cells.Item(26, 38).NumberFormat = "@"
$ws.Cells.Item(26, 38).Value = 'Si'
$ws.Cells.Item(26, 38).NumberFormat = "General"
$ws.Cells.Item(26, 39).NumberFormat = "@"
$ws.Cells.Item(26, 39).Value = 'es tracta del marhe de la linea de transecció previa'
$ws.Cells.Item(26, 39).NumberFormat = "General"
$ws.Cells.Item(26, 40).Value = 44522
$ws.Cells.Item(26, 41).NumberFormat = "@"
$ws.Cells.Item(26, 41).Value = 'No'
$ws.Cells.Item(26, 41).NumberFormat = "General"
$ws.Cells.Item(26, 42).NumberFormat = "@"
$ws.Cells.Item(26, 42).Value = 'No'
$ws.Cells.Item(26, 42).NumberFormat = "General"
$ws.Cells.Item(26, 43).NumberFormat = "@"
$ws.Cells.Item(26, 43).Value = 'Viu'
$ws.Cells.Item(26, 43).NumberFormat = "General"
$ws.Cells.Item(26, 44).NumberFormat = "@"
$ws.Cells.Item(26, 44).Value = 'Si'
$ws.Cells.Item(26, 44).NumberFormat = "General"
$ws.Cells.Item(26, 45).NumberFormat = "@"
$ws.Cells.Item(26, 45).Value = 'Si'
$ws.Cells.Item(26, 45).NumberFormat = "General"
$ws.Cells.Item(26, 46).NumberFormat = "@"
$ws.Cells.Item(26, 46).Value = 'No'
$ws.Cells.Item(26, 46).NumberFormat = "General"
$ws.Cells.Item(26, 47).NumberFormat = "@"
$ws.Cells.Item(26, 47).Value = 'No'
$ws.Cells.Item(26, 47).NumberFormat = "General"
$ws.Cells.Item(26, 48).NumberFormat = "@"
$ws.Cells.Item(26, 48).Value = 'No'
$ws.Cells.Item(26, 48).NumberFormat = "General"
$ws.Cells.Item(26, 49).NumberFormat = "@"
$ws.Cells.Item(26, 49).Value = 'No'
$ws.Cells.Item(26, 49).NumberFormat = "General"
$ws.Cells.Item(26, 50).NumberFormat = "@"
$ws.Cells.Item(26, 50).Value = 'Si'
$ws.Cells.Item(26, 50).NumberFormat = "General"
$ws.Cells.Item(26, 51).NumberFormat = "@"
$ws.Cells.Item(26, 51).Value = 'Si'
$ws.Cells.Item(26, 51).NumberFormat = "General"
$ws.Cells.Item(26, 52).NumberFormat = "@"
$ws.Cells.Item(26, 52).Value = 'Si'
$ws.Cells.Item(26, 52).NumberFormat = "General"
$ws.Cells.Item(26, 53).Value = 6
$ws.Cells.Item(26, 54).Value = 44610.47984626157
$ws.Cells.Item(26, 55).NumberFormat = "@"
$ws.Cells.Item(26, 55).Value = 'Falta alguna variable, revisar'
$ws.Cells.Item(26, 55).NumberFormat = "General"
$ws.Cells.Item(26, 56).Value = 43342
$ws.Cells.Item(26, 57).NumberFormat = "@"
$ws.Cells.Item(26, 57).Value = 'No'
$ws.Cells.Item(26, 57).NumberFormat = "General"
$ws.Cells.Item(26, 59).Value = 1
$ws.Cells.Item(26, 63).NumberFormat = "@"
$ws.Cells.Item(26, 63).Value = 'Segmentectomia1a8'
$ws.Cells.Item(26, 63).NumberFormat = "General"

# --- Row 27 ---
$ws.Cells.Item(27, 1).NumberFormat = "@"
$ws.Cells.Item(27, 1).Value = 'Hepatectomia dreta'
$ws.Cells.Item(27, 1).NumberFormat = "General"
$ws.Cells.Item(27, 2).Value = 1695
$ws.Cells.Item(27, 3).NumberFormat = "@"
$ws.Cells.Item(27, 3).Value = '09/07/2018'
$ws.Cells.Item(27, 3).NumberFormat = "General"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '11/02/2019'
$ws.Cells.Item(27, 4).NumberFormat = "General"
$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = '25/06/2018'
$ws.Cells.Item(27, 5).NumberFormat = "General"
$ws.Cells.Item(27, 6).NumberFormat = "@"
$ws.Cells.Item(27, 6).Value = 'VII'
$ws.Cells.Item(27, 6).NumberFormat = "General"
$ws.Cells.Item(27, 8).NumberFormat = "@"
$ws.Cells.Item(27, 8).Value = 'Jordi'
$ws.Cells.Item(27, 8).NumberFormat = "General"
$ws.Cells.Item(27, 9).NumberFormat = "@"
$ws.Cells.Item(27, 9).Value = 'Morillas2'
$ws.Cells.Item(27, 9).NumberFormat = "General"
$ws.Cells.Item(27, 10).NumberFormat = "@"
$ws.Cells.Item(27, 10).Value = 'Esteban2'
$ws.Cells.Item(27, 10).NumberFormat = "General"
$ws.Cells.Item(27, 11).Value = 13296015
$ws.Cells.Item(27, 12).NumberFormat = "@"
$ws.Cells.Item(27, 12).Value = 'No'
$ws.Cells.Item(27, 12).NumberFormat = "General"
$ws.Cells.Item(27, 15).NumberFormat = "@"
$ws.Cells.Item(27, 15).Value = '79'
$ws.Cells.Item(27, 15).NumberFormat = "General"
$ws.Cells.Item(27, 16).Value = 178
$ws.Cells.Item(27, 17).Value = 25
$ws.Cells.Item(27, 18).Value = 3
$ws.Cells.Item(27, 19).NumberFormat = "@"
$ws.Cells.Item(27, 19).Value = 'No'
$ws.Cells.Item(27, 19).NumberFormat = "General"
$ws.Cells.Item(27, 20).Value = 43858
$ws.Cells.Item(27, 21).NumberFormat = "@"
$ws.Cells.Item(27, 21).Value = 'Resecció Major (>= 3 segm)'
$ws.Cells.Item(27, 21).NumberFormat = "General"
$ws.Cells.Item(27, 22).NumberFormat = "@"
$ws.Cells.Item(27, 22).Value = 'hepatectomia dreta'
$ws.Cells.Item(27, 22).NumberFormat = "General"
$ws.Cells.Item(27, 23).NumberFormat = "@"
$ws.Cells.Item(27, 23).Value = 'Oberta'
$ws.Cells.Item(27, 23).NumberFormat = "General"
$ws.Cells.Item(27, 24).NumberFormat = "@"
$ws.Cells.Item(27, 24).Value = 'Si, com a primer temps quirúrgic'
$ws.Cells.Item(27, 24).NumberFormat = "General"
$ws.Cells.Item(27, 25).NumberFormat = "@"
$ws.Cells.Item(27, 25).Value = 'No'
$ws.Cells.Item(27, 25).NumberFormat = "General"
$ws.Cells.Item(27, 26).NumberFormat = "@"
$ws.Cells.Item(27, 26).Value = 'Impressió R1'
$ws.Cells.Item(27, 26).NumberFormat = "General"
$ws.Cells.Item(27, 27).Value = 1
$ws.Cells.Item(27, 28).Value = 3
$ws.Cells.Item(27, 29).NumberFormat = "@"
$ws.Cells.Item(27, 29).Value = 'No'
$ws.Cells.Item(27, 29).NumberFormat = "General"
$ws.Cells.Item(27, 30).NumberFormat = "@"
$ws.Cells.Item(27, 30).Value = 'Si'
$ws.Cells.Item(27, 30).NumberFormat = "General"
$ws.Cells.Item(27, 31).NumberFormat = "@"
$ws.Cells.Item(27, 31).Value = 'No'
$ws.Cells.Item(27, 31).NumberFormat = "General"
$ws.Cells.Item(27, 32).NumberFormat = "@"
$ws.Cells.Item(27, 32).Value = 'Si'
$ws.Cells.Item(27, 32).NumberFormat = "General"
$ws.Cells.Item(27, 33).NumberFormat = "@"
$ws.Cells.Item(27, 33).Value = 'IIIb'
$ws.Cells.Item(27, 33).NumberFormat = "General"
$ws.Cells.Item(27, 34).Value = 61
$ws.Cells.Item(27, 35).Value = 1
$ws.Cells.Item(27, 36).Value = 3
$ws.Cells.Item(27, 37).Value = 0
$ws.Cells.Item(27, 38).NumberFormat = "@"
$ws.Cells.Item(27, 38).Value = 'Si'
$ws.Cells.Item(27, 38).NumberFormat = "General"
$ws.Cells.Item(27, 39).NumberFormat = "@"
$ws.Cells.Item(27, 39).Value = 'ampliacio quirurgica'
$ws.Cells.Item(27, 39).NumberFormat = "General"
$ws.Cells.Item(27, 40).Value = 43983
$ws.Cells.Item(27, 41).NumberFormat = "@"
$ws.Cells.Item(27, 41).Value = 'No'
$ws.Cells.Item(27, 41).NumberFormat = "General"
$ws.Cells.Item(27, 42).NumberFormat = "@"
$ws.Cells.Item(27, 42).Value = 'No'
$ws.Cells.Item(27, 42).NumberFormat = "General"
$ws.Cells.Item(27, 43).NumberFormat = "@"
$ws.Cells.Item(27, 43).Value = 'Viu'
$ws.Cells.Item(27, 43).NumberFormat = "General"
$ws.Cells.Item(27, 44).NumberFormat = "@"
$ws.Cells.Item(27, 44).Value = 'Si'
$ws.Cells.Item(27, 44).NumberFormat = "General"
$ws.Cells.Item(27, 45).NumberFormat = "@"
$ws.Cells.Item(27, 45).Value = 'Si'
$ws.Cells.Item(27, 45).NumberFormat = "General"
$ws.Cells.Item(27, 46).NumberFormat = "@"
$ws.Cells.Item(27, 46).Value = 'No'
$ws.Cells.Item(27, 46).NumberFormat = "General"
$ws.Cells.Item(27, 47).NumberFormat = "@"
$ws.Cells.Item(27, 47).Value = 'No'
$ws.Cells.Item(27, 47).NumberFormat = "General"
$ws.Cells.Item(27, 48).NumberFormat = "@"
$ws.Cells.Item(27, 48).Value = 'Si'
$ws.Cells.Item(27, 48).NumberFormat = "General"
$ws.Cells.Item(27, 49).NumberFormat = "@"
$ws.Cells.Item(27, 49).Value = 'No'
$ws.Cells.Item(27, 49).NumberFormat = "General"
$ws.Cells.Item(27, 50).NumberFormat = "@"
$ws.Cells.Item(27, 50).Value = 'No'
$ws.Cells.Item(27, 50).NumberFormat = "General"
$ws.Cells.Item(27, 51).NumberFormat = "@"
$ws.Cells.Item(27, 51).Value = 'No'
$ws.Cells.Item(27, 51).NumberFormat = "General"
$ws.Cells.Item(27, 52).NumberFormat = "@"
$ws.Cells.Item(27, 52).Value = 'No'
$ws.Cells.Item(27, 52).NumberFormat = "General"
$ws.Cells.Item(27, 53).Value = 55
$ws.Cells.Item(27, 54).Value = 44610.50603524306
$ws.Cells.Item(27, 55).NumberFormat = "@"
$ws.Cells.Item(27, 55).Value = 'Error en alguna fechas formateadas, revisar'
$ws.Cells.Item(27, 55).NumberFormat = "General"
$ws.Cells.Item(27, 56).Value = 43446
$ws.Cells.Item(27, 57).NumberFormat = "@"
$ws.Cells.Item(27, 57).Value = 'No'
$ws.Cells.Item(27, 57).NumberFormat = "General"
$ws.Cells.Item(27, 58).NumberFormat = "@"
$ws.Cells.Item(27, 58).Value = 'oclusió, peritonitis fecaoidea'
$ws.Cells.Item(27, 58).NumberFormat = "General"
$ws.Cells.Item(27, 59).Value = 2
$ws.Cells.Item(27, 63).NumberFormat = "@"
$ws.Cells.Item(27, 63).Value = 'hepatectomiaDreta'
$ws.Cells.Item(27, 63).NumberFormat = "General"

# --- Row 28 ---
$ws.Cells.Item(28, 1).NumberFormat = "@"
$ws.Cells.Item(28, 1).Value = 'Hepatectomia dreta'
$ws.Cells.Item(28, 1).NumberFormat = "General"
$ws.Cells.Item(28, 2).Value = 1695
$ws.Cells.Item(28, 3).NumberFormat = "@"
$ws.Cells.Item(28, 3).Value = '09/07/2018'
$ws.Cells.Item(28, 3).NumberFormat = "General"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '11/02/2019'
$ws.Cells.Item(28, 4).NumberFormat = "General"
$ws.Cells.Item(28, 5).NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = '25/06/2018'
$ws.Cells.Item(28, 5).NumberFormat = "General"
$ws.Cells.Item(28, 6).NumberFormat = "@"
$ws.Cells.Item(28, 6).Value = 'VII'
$ws.Cells.Item(28, 6).NumberFormat = "General"
$ws.Cells.Item(28, 8).NumberFormat = "@"
$ws.Cells.Item(28, 8).Value = 'Jordi'
$ws.Cells.Item(28, 8).NumberFormat = "General"
$ws.Cells.Item(28, 9).NumberFormat = "@"
$ws.Cells.Item(28, 9).Value = 'Morillas2'
$ws.Cells.Item(28, 9).NumberFormat = "General"
$ws.Cells.Item(28, 10).NumberFormat = "@"
$ws.Cells.Item(28, 10).Value = 'Esteban2'
$ws.Cells.Item(28, 10).NumberFormat = "General"
$ws.Cells.Item(28, 11).Value = 13296015
$ws.Cells.Item(28, 12).NumberFormat = "@"
$ws.Cells.Item(28, 12).Value = 'No'
$ws.Cells.Item(28, 12).NumberFormat = "General"
$ws.Cells.Item(28, 15).NumberFormat = "@"
$ws.Cells.Item(28, 15).Value = '79'
$ws.Cells.Item(28, 15).NumberFormat = "General"
$ws.Cells.Item(28, 16).Value = 178
$ws.Cells.Item(28, 17).Value = 25
$ws.Cells.Item(28, 18).Value = 3
$ws.Cells.Item(28, 19).NumberFormat = "@"
$ws.Cells.Item(28, 19).Value = 'No'
$ws.Cells.Item(28, 19).NumberFormat = "General"
$ws.Cells.Item(28, 20).Value = 43858
$ws.Cells.Item(28, 21).NumberFormat = "@"
$ws.Cells.Item(28, 21).Value = 'Resecció Major (>= 3 segm)'
$ws.Cells.Item(28, 21).NumberFormat = "General"
$ws.Cells.Item(28, 22).NumberFormat = "@"
$ws.Cells.Item(28, 22).Value = 'hepatectomia dreta'
$ws.Cells.Item(28, 22).NumberFormat = "General"
$ws.Cells.Item(28, 23).NumberFormat = "@"
$ws.Cells.Item(28, 23).Value = 'Oberta'
$ws.Cells.Item(28, 23).NumberFormat = "General"
$ws.Cells.Item(28, 24).NumberFormat = "@"
$ws.Cells.Item(28, 24).Value = 'Si, com a primer temps quirúrgic'
$ws.Cells.Item(28, 24).NumberFormat = "General"
$ws.Cells.Item(28, 25).NumberFormat = "@"
$ws.Cells.Item(28, 25).Value = 'No'
$ws.Cells.Item(28, 25).NumberFormat = "General"
$ws.Cells.Item(28, 26).NumberFormat = "@"
$ws.Cells.Item(28, 26).Value = 'Impressió R1'
$ws.Cells.Item(28, 26).NumberFormat = "General"
$ws.Cells.Item(28, 27).Value = 1
$ws.Cells.Item(28, 28).Value = 3
$ws.Cells.Item(28, 29).NumberFormat = "@"
$ws.Cells.Item(28, 29).Value = 'No'
$ws.Cells.Item(28, 29).NumberFormat = "General"
$ws.Cells.Item(28, 30).NumberFormat = "@"
$ws.Cells.Item(28, 30).Value = 'Si'
$ws.Cells.Item(28, 30).NumberFormat = "General"
$ws.Cells.Item(28, 31).NumberFormat = "@"
$ws.Cells.Item(28, 31).Value = 'No'
$ws.Cells.Item(28, 31).NumberFormat = "General"
$ws.Cells.Item(28, 32).NumberFormat = "@"
$ws.Cells.Item(28, 32).Value = 'Si'
$ws.Cells.Item(28, 32).NumberFormat = "General"
$ws.Cells.Item(28, 33).NumberFormat = "@"
$ws.Cells.Item(28, 33).Value = 'IIIb'
$ws.Cells.Item(28, 33).NumberFormat = "General"
$ws.Cells.Item(28, 34).Value = 61
$ws.Cells.Item(28, 35).Value = 1
$ws.Cells.Item(28, 36).Value = 3
$ws.Cells.Item(28, 37).Value = 0
$ws.Cells.Item(28, 38).NumberFormat = "@"
$ws.Cells.Item(28, 38).Value = 'Si'
$ws.Cells.Item(28, 38).NumberFormat = "General"
$ws.Cells.Item(28, 39).NumberFormat = "@"
$ws.Cells.Item(28, 39).Value = 'ampliacio quirurgica'
$ws.Cells.Item(28, 39).NumberFormat = "General"
$ws.Cells.Item(28, 40).Value = 43983
$ws.Cells.Item(28, 41).NumberFormat = "@"
$ws.Cells.Item(28, 41).Value = 'No'
$ws.Cells.Item(28, 41).NumberFormat = "General"
$ws.Cells.Item(28, 42).NumberFormat = "@"
$ws.Cells.Item(28, 42).Value = 'No'
$ws.Cells.Item(28, 42).NumberFormat = "General"
$ws.Cells.Item(28, 43).NumberFormat = "@"
$ws.Cells.Item(28, 43).Value = 'Viu'
$ws.Cells.Item(28, 43).NumberFormat = "General"
$ws.Cells.Item(28, 44).NumberFormat = "@"
$ws.Cells.Item(28, 44).Value = 'Si'
$ws.Cells.Item(28, 44).NumberFormat = "General"
$ws.Cells.Item(28, 45).NumberFormat = "@"
$ws.Cells.Item(28, 45).Value = 'Si'
$ws.Cells.Item(28, 45).NumberFormat = "General"
$ws.Cells.Item(28, 46).NumberFormat = "@"
$ws.Cells.Item(28, 46).Value = 'No'
$ws.Cells.Item(28, 46).NumberFormat = "General"
$ws.Cells.Item(28, 47).NumberFormat = "@"
$ws.Cells.Item(28, 47).Value = 'No'
$ws.Cells.Item(28, 47).NumberFormat = "General"
$ws.Cells.Item(28, 48).NumberFormat = "@"
$ws.Cells.Item(28, 48).Value = 'Si'
$ws.Cells.Item(28, 48).NumberFormat = "General"
$ws.Cells.Item(28, 49).NumberFormat = "@"
$ws.Cells.Item(28, 49).Value = 'No'
$ws.Cells.Item(28, 49).NumberFormat = "General"
$ws.Cells.Item(28, 50).NumberFormat = "@"
$ws.Cells.Item(28, 50).Value = 'No'
$ws.Cells.Item(28, 50).NumberFormat = "General"
$ws.Cells.Item(28, 51).NumberFormat = "@"
$ws.Cells.Item(28, 51).Value = 'No'
$ws.Cells.Item(28, 51).NumberFormat = "General"
$ws.Cells.Item(28, 52).NumberFormat = "@"
$ws.Cells.Item(28, 52).Value = 'No'
$ws.Cells.Item(28, 52).NumberFormat = "General"
$ws.Cells.Item(28, 53).Value = 55
$ws.Cells.Item(28, 54).Value = 44610.50612837963
$ws.Cells.Item(28, 55).NumberFormat = "@"
$ws.Cells.Item(28, 55).Value = 'Error en alguna fechas formateadas, revisar'
$ws.Cells.Item(28, 55).NumberFormat = "General"
$ws.Cells.Item(28, 56).Value = 43446
$ws.Cells.Item(28, 57).NumberFormat = "@"
$ws.Cells.Item(28, 57).Value = 'No'
$ws.Cells.Item(28, 57).NumberFormat = "General"
$ws.Cells.Item(28, 58).NumberFormat = "@"
$ws.Cells.Item(28, 58).Value = 'oclusió, peritonitis fecaoidea'
$ws.Cells.Item(28, 58).NumberFormat = "General"
$ws.Cells.Item(28, 59).Value = 2
$ws.Cells.Item(28, 63).NumberFormat = "@"
$ws.Cells.Item(28, 63).Value = 'hepatectomiaDreta'
$ws.Cells.Item(28, 63).NumberFormat = "General"

# --- Row 29 ---
$ws.Cells.Item(29, 1).NumberFormat = "@"
$ws.Cells.Item(29, 1).Value = 'Hepatectomia dreta'
$ws.Cells.Item(29, 1).NumberFormat = "General"
$ws.Cells.Item(29, 2).Value = 1695
$ws.Cells.Item(29, 3).NumberFormat = "@"
$ws.Cells.Item(29, 3).Value = '09/07/2018'
$ws.Cells.Item(29, 3).NumberFormat = "General"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '11/02/2019'
$ws.Cells.Item(29, 4).NumberFormat = "General"
$ws.Cells.Item(29, 5).NumberFormat = "@"
$ws.Cells.Item(29, 5).Value = '25/06/2018'
$ws.Cells.Item(29, 5).NumberFormat = "General"
$ws.Cells.Item(29, 6).NumberFormat = "@"
$ws.Cells.Item(29, 6).Value = 'VII'
$ws.Cells.Item(29, 6).NumberFormat = "General"
$ws.Cells.Item(29, 8).NumberFormat = "@"
$ws.Cells.Item(29, 8).Value = 'Jordi'
$ws.Cells.Item(29, 8).NumberFormat = "General"
$ws.Cells.Item(29, 9).NumberFormat = "@"
$ws.Cells.Item(29, 9).Value = 'Morillas2'
$ws.Cells.Item(29, 9).NumberFormat = "General"
$ws.Cells.Item(29, 10).NumberFormat = "@"
$ws.Cells.Item(29, 10).Value = 'Esteban2'
$ws.Cells.Item(29, 10).NumberFormat = "General"
$ws.Cells.Item(29, 11).Value = 13296015
$ws.Cells.Item(29, 12).NumberFormat = "@"
$ws.Cells.Item(29, 12).Value = 'No'
$ws.Cells.Item(29, 12).NumberFormat = "General"
$ws.Cells.Item(29, 15).NumberFormat = "@"
$ws.Cells.Item(29, 15).Value = '79'
$ws.Cells.Item(29, 15).NumberFormat = "General"
$ws.Cells.Item(29, 16).Value = 178
$ws.Cells.Item(29, 17).Value = 25
$ws.Cells.Item(29, 18).Value = 3
$ws.Cells.Item(29, 19).NumberFormat = "@"
$ws.Cells.Item(29, 19).Value = 'No'
$ws.Cells.Item(29, 19).NumberFormat = "General"
$ws.Cells.Item(29, 20).Value = 43858
$ws.Cells.Item(29, 21).NumberFormat = "@"
$ws.Cells.Item(29, 21).Value = 'Resecció Major (>= 3 segm)'
$ws.Cells.Item(29, 21).NumberFormat = "General"
$ws.Cells.Item(29, 22).NumberFormat = "@"
$ws.Cells.Item(29, 22).Value = 'hepatectomia dreta'
$ws.Cells.Item(29, 22).NumberFormat = "General"
$ws.Cells.Item(29, 23).NumberFormat = "@"
$ws.Cells.Item(29, 23).Value = 'Oberta'
$ws.Cells.Item(29, 23).NumberFormat = "General"
$ws.Cells.Item(29, 24).NumberFormat = "@"
$ws.Cells.Item(29, 24).Value = 'Si, com a primer temps quirúrgic'
$ws.Cells.Item(29, 24).NumberFormat = "General"
$ws.Cells.Item(29, 25).NumberFormat = "@"
$ws.Cells.Item(29, 25).Value = 'No'
$ws.Cells.Item(29, 25).NumberFormat = "General"
$ws.Cells.Item(29, 26).NumberFormat = "@"
$ws.Cells.Item(29, 26).Value = 'Impressió R1'
$ws.Cells.Item(29, 26).NumberFormat = "General"
$ws.Cells.Item(29, 27).Value = 1
$ws.Cells.Item(29, 28).Value = 3
$ws.Cells.Item(29, 29).NumberFormat = "@"
$ws.Cells.Item(29, 29).Value = 'No'
$ws.Cells.Item(29, 29).NumberFormat = "General"
$ws.Cells.Item(29, 30).NumberFormat = "@"
$ws.Cells.Item(29, 30).Value = 'Si'
$ws.Cells.Item(29, 30).NumberFormat = "General"
$ws.Cells.Item(29, 31).NumberFormat = "@"
$ws.Cells.Item(29, 31).Value = 'No'
$ws.Cells.Item(29, 31).NumberFormat = "General"
$ws.Cells.Item(29, 32).NumberFormat = "@"
$ws.Cells.Item(29, 32).Value = 'Si'
$ws.Cells.Item(29, 32).NumberFormat = "General"
$ws.Cells.Item(29, 33).NumberFormat = "@"
$ws.Cells.Item(29, 33).Value = 'IIIb'
$ws.Cells.Item(29, 33).NumberFormat = "General"
$ws.Cells.Item(29, 34).Value = 61
$ws.Cells.Item(29, 35).Value = 1
$ws.Cells.Item(29, 36).Value = 3
$ws.Cells.Item(29, 37).Value = 0
$ws.Cells.Item(29, 38).NumberFormat = "@"
$ws.Cells.Item(29, 38).Value = 'Si'
$ws.Cells.Item(29, 38).NumberFormat = "General"
$ws.Cells.Item(29, 39).NumberFormat = "@"
$ws.Cells.Item(29, 39).Value = 'ampliacio quirurgica'
$ws.Cells.Item(29, 39).NumberFormat = "General"
$ws.Cells.Item(29, 40).Value = 43983
$ws.Cells.Item(29, 41).NumberFormat = "@"
$ws.Cells.Item(29, 41).Value = 'No'
$ws.Cells.Item(29, 41).NumberFormat = "General"
$ws.Cells.Item(29, 42).NumberFormat = "@"
$ws.Cells.Item(29, 42).Value = 'No'
$ws.Cells.Item(29, 42).NumberFormat = "General"
$ws.Cells.Item(29, 43).NumberFormat = "@"
$ws.Cells.Item(29, 43).Value = 'Viu'
$ws.Cells.Item(29, 43).NumberFormat = "General"
$ws.Cells.Item(29, 44).NumberFormat = "@"
$ws.Cells.Item(29, 44).Value = 'Si'
$ws.Cells.Item(29, 44).NumberFormat = "General"
$ws.Cells.Item(29, 45).NumberFormat = "@"
$ws.Cells.Item(29, 45).Value = 'Si'
$ws.Cells.Item(29, 45).NumberFormat = "General"
$ws.Cells.Item(29, 46).NumberFormat = "@"
$ws.Cells.Item(29, 46).Value = 'No'
$ws.Cells.Item(29, 46).NumberFormat = "General"
$ws.Cells.Item(29, 47).NumberFormat = "@"
$ws.Cells.Item(29, 47).Value = 'No'
$ws.Cells.Item(29, 47).NumberFormat = "General"
$ws.Cells.Item(29, 48).NumberFormat = "@"
$ws.Cells.Item(29, 48).Value = 'Si'
$ws.Cells.Item(29, 48).NumberFormat = "General"
$ws.Cells.Item(29, 49).NumberFormat = "@"
$ws.Cells.Item(29, 49).Value = 'No'
$ws.Cells.Item(29, 49).NumberFormat = "General"
$ws.Cells.Item(29, 50).NumberFormat = "@"
$ws.Cells.Item(29, 50).Value = 'No'
$ws.Cells.Item(29, 50).NumberFormat = "General"
$ws.Cells.Item(29, 51).NumberFormat = "@"
$ws.Cells.Item(29, 51).Value = 'No'
$ws.Cells.Item(29, 51).NumberFormat = "General"
$ws.Cells.Item(29, 52).NumberFormat = "@"
$ws.Cells.Item(29, 52).Value = 'No'
$ws.Cells.Item(29, 52).NumberFormat = "General"
$ws.Cells.Item(29, 53).Value = 55
$ws.Cells.Item(29, 54).Value = 44610.50628342592
$ws.Cells.Item(29, 55).NumberFormat = "@"
$ws.Cells.Item(29, 55).Value = 'Error en alguna fechas formateadas, revisar'
$ws.Cells.Item(29, 55).NumberFormat = "General"
$ws.Cells.Item(29, 56).Value = 43446
$ws.Cells.Item(29, 57).NumberFormat = "@"
$ws.Cells.Item(29, 57).Value = 'No'
$ws.Cells.Item(29, 57).NumberFormat = "General"
$ws.Cells.Item(29, 58).NumberFormat = "@"
$ws.Cells.Item(29, 58).Value = 'oclusió, peritonitis fecaoidea'
$ws.Cells.Item(29, 58).NumberFormat = "General"
$ws.Cells.Item(29, 59).Value = 2
$ws.Cells.Item(29, 63).NumberFormat = "@"
$ws.Cells.Item(29, 63).Value = 'hepatectomiaDreta'
$ws.Cells.Item(29, 63).NumberFormat = "General"

Write-Host "edit complete"
